# Updates the cryptos price/volume snapshot table on Sheet1 (row 2 = Bitcoin ...
# row 51 = last coin) to match the latest scrape.
#
# Columns: B=Coin name, C=Link, D=Price (plain text, not numeric - values such
# as "27.851.37" are thousands-grouped prices stored as text, and plain
# decimals such as "336.90" must stay text too, so they are written with a
# leading apostrophe to stop Excel from auto-converting them to numbers),
# E=Volume(1h) change percentage (plain text, kept as "  +0.15%  " style
# strings with surrounding padding spaces).
#
# Row 51 switched from PancakeSwap to EOS in this update, so its coin name,
# link, price and percentage all change together.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.851.37'
$ws.Range('E2').Value = '  -0.07%  '
$ws.Range('D3').Value = '1.884.79'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').Value = '''336.90'
$ws.Range('E5').Value = '  +0.69%  '
$ws.Range('D6').Value = '''1.011'
$ws.Range('E6').Value = '  +0.27%  '
$ws.Range('D7').Value = '''0.4674'
$ws.Range('E7').Value = '  -1.36%  '
$ws.Range('D8').Value = '''0.3958'
$ws.Range('E8').Value = '  +0.60%  '
$ws.Range('D9').Value = '''46.05'
$ws.Range('E9').Value = '  -3.34%  '
$ws.Range('D10').Value = '''0.08020'
$ws.Range('E10').Value = '  -0.81%  '
$ws.Range('D11').Value = '''1.012'
$ws.Range('E11').Value = '  -1.57%  '
$ws.Range('D12').Value = '''22.00'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '1.893.29'
$ws.Range('E13').Value = '  +0.33%  '
$ws.Range('D14').Value = '''6.005'
$ws.Range('E14').Value = '  +0.18%  '
$ws.Range('D15').Value = '''7.285'
$ws.Range('E15').Value = '  +1.76%  '
$ws.Range('D16').Value = '''1.014'
$ws.Range('D17').Value = '''89.30'
$ws.Range('E17').Value = '  +2.12%  '
$ws.Range('D18').Value = '''0.06717'
$ws.Range('E18').Value = '  -0.48%  '
$ws.Range('D19').Value = '''0.00001048'
$ws.Range('E19').Value = '  -0.43%  '
$ws.Range('D20').Value = '''17.37'
$ws.Range('E20').Value = '  +0.05%  '
$ws.Range('D21').Value = '''1.011'
$ws.Range('E21').Value = '  +0.29%  '
$ws.Range('D22').Value = '27.878.32'
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('D23').Value = '''5.504'
$ws.Range('E23').Value = '  -0.44%  '
$ws.Range('D24').Value = '''11.02'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('D25').Value = '''2.315'
$ws.Range('E25').Value = '  -0.73%  '
$ws.Range('D26').Value = '2.112.72'
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '''159.15'
$ws.Range('E27').Value = '  +0.01%  '
$ws.Range('D28').Value = '''19.85'
$ws.Range('E28').Value = '  -1.28%  '
$ws.Range('D29').Value = '''2.159'
$ws.Range('E29').Value = '  +2.31%  '
$ws.Range('D30').Value = '''5.501'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').Value = '''122.07'
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('D32').Value = '''0.9860'
$ws.Range('E32').Value = '  +0.66%  '
$ws.Range('D33').Value = '''0.09468'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('D34').Value = '''3.643'
$ws.Range('E34').Value = '  +0.46%  '
$ws.Range('D35').Value = '''5.346'
$ws.Range('E35').Value = '  -0.37%  '
$ws.Range('D36').Value = '''1.357'
$ws.Range('E36').Value = '  -6.66%  '
$ws.Range('D37').Value = '''0.06088'
$ws.Range('E37').Value = '  -1.42%  '
$ws.Range('D38').Value = '''0.02247'
$ws.Range('E38').Value = '  -1.22%  '
$ws.Range('D39').Value = '''8.385'
$ws.Range('E39').Value = '  +3.56%  '
$ws.Range('D40').Value = '''1.202'
$ws.Range('E40').Value = '  -1.56%  '
$ws.Range('D41').Value = '''1.010'
$ws.Range('E41').Value = '  +0.22%  '
$ws.Range('D42').Value = '''0.5998'
$ws.Range('E42').Value = '  -0.34%  '
$ws.Range('D43').Value = '''0.1898'
$ws.Range('E43').Value = '  -0.14%  '
$ws.Range('D44').Value = '''10.41'
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('D45').Value = '''1.248'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('D46').Value = '''0.5669'
$ws.Range('E46').Value = '  -0.88%  '
$ws.Range('D47').Value = '''12.34'
$ws.Range('E47').Value = '  +0.34%  '
$ws.Range('D48').Value = '''1.946'
$ws.Range('E48').Value = '  +0.12%  '
$ws.Range('D49').Value = '''0.06789'
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('D50').Value = '''112.79'
$ws.Range('E50').Value = '  -0.72%  '
$ws.Range('B51').Value = 'EOS'
$ws.Range('C51').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D51').Value = '''1.066'
$ws.Range('E51').Value = '  -0.65%  '
